# PCB and DOCS update: Page structure (WIP) and BOM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# B3: stabilizátor 3V3 count 2 -> 0
$ws.Range("B3").Value = 0

# B4: rezistor 10k count 0 -> 5
$ws.Range("B4").Value = 5

# B23: Bórykův spínač count 3 -> 1
$ws.Range("B23").Value = 1

# Row 32 was blank; add new BOM line "vložky M3 " with quantity 30
$ws.Range("A32").Value = "vložky M3 "
$ws.Range("B32").Value = 30
